$wb = $excel.ActiveWorkbook

$wsTutor = $wb.Worksheets.Item("Tutor")
# Remove the "Highly Rated Mentor" / "Have a rating of 5 in 3 courses" row (row 22)
$wsTutor.Rows.Item(22).Delete()
# Remove the "Top-Rated Instructor" / "Have a rating of 5 in a course" row (row 21)
$wsTutor.Rows.Item(21).Delete()
# Remove the "Beginner Tutor" / "Have a role to be tutor" row (row 1)
$wsTutor.Rows.Item(1).Delete()

$wsStudent = $wb.Worksheets.Item("Student")
# Remove the "Beginner Student" / "Have a role to be student" row (row 1)
$wsStudent.Rows.Item(1).Delete()

$wsTutor.Range("F9").Select() | Out-Null
$wsStudent.Range("J11").Select() | Out-Null
